$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 0.6113651253405055
$ws.Range("J2").Value = 0.6113651253405055
$ws.Range("M2").Value = 173.5452066666667
$ws.Range("N2").Value = 520.63562
$ws.Range("O2").Value = 0.6098887991422922
$ws.Range("P2").Value = 0.6098887991422922
$ws.Range("Q2").Value = 5.128318705402222
$ws.Range("R2").Value = 46.15486834862
$ws.Range("S2").Value = 0.3728647421313979
$ws.Range("T2").Value = 0.3728647421313979
$ws.Range("I3").Value = 0.6113651253405055
$ws.Range("J3").Value = 0.6113651253405055
$ws.Range("O3").Value = 0.003264284357140855
$ws.Range("P3").Value = 0.003264284357140855
$ws.Range("S3").Value = 0.00199566961515047
$ws.Range("T3").Value = 0.00199566961515047
$ws.Range("I4").Value = 0.6113651253405055
$ws.Range("J4").Value = 0.6113651253405055
$ws.Range("M4").Value = 54.64271666666667
$ws.Range("N4").Value = 163.92815
$ws.Range("O4").Value = 0.192030546333187
$ws.Range("P4").Value = 0.192030546333187
$ws.Range("Q4").Value = 1.614710491738889
$ws.Range("R4").Value = 14.53239442565
$ws.Range("S4").Value = 0.1174007790281946
$ws.Range("T4").Value = 0.1174007790281946
$ws.Range("I5").Value = 0.6113651253405055
$ws.Range("J5").Value = 0.6113651253405055
$ws.Range("M5").Value = 1.069012
$ws.Range("N5").Value = 3.207036
$ws.Range("O5").Value = 0.00375682196858928
$ws.Range("P5").Value = 0.00375682196858928
$ws.Range("Q5").Value = 0.03158966093733333
$ws.Range("R5").Value = 0.284306948436
$ws.Range("S5").Value = 0.002296789933708549
$ws.Range("T5").Value = 0.00229678993370855
$ws.Range("I6").Value = 0.6113651253405055
$ws.Range("J6").Value = 0.6113651253405055
$ws.Range("M6").Value = 54.36641700000001
$ws.Range("N6").Value = 163.099251
$ws.Range("O6").Value = 0.1910595481987908
$ws.Range("P6").Value = 0.1910595481987908
$ws.Range("Q6").Value = 1.606545744489
$ws.Range("R6").Value = 14.458911700401
$ws.Range("S6").Value = 0.1168071446320541
$ws.Range("T6").Value = 0.1168071446320541
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01878466666666667
$ws.Range("H7").Value = 0.056354
$ws.Range("I7").Value = 0.3886348746594945
$ws.Range("J7").Value = 0.3886348746594945
$ws.Range("M7").Value = 173.5452066666667
$ws.Range("N7").Value = 520.63562
$ws.Range("O7").Value = 0.6098887991422922
$ws.Range("P7").Value = 0.6098887991422922
$ws.Range("Q7").Value = 3.259988858831111
$ws.Range("R7").Value = 29.33989972948
$ws.Range("S7").Value = 0.2370240570108944
$ws.Range("T7").Value = 0.2370240570108944
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.01878466666666667
$ws.Range("H8").Value = 0.056354
$ws.Range("I8").Value = 0.3886348746594945
$ws.Range("J8").Value = 0.3886348746594945
$ws.Range("O8").Value = 0.003264284357140855
$ws.Range("P8").Value = 0.003264284357140855
$ws.Range("Q8").Value = 0.01744831295688889
$ws.Range("R8").Value = 0.157034816612
$ws.Range("S8").Value = 0.001268614741990385
$ws.Range("T8").Value = 0.001268614741990385
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.01878466666666667
$ws.Range("H9").Value = 0.056354
$ws.Range("I9").Value = 0.3886348746594945
$ws.Range("J9").Value = 0.3886348746594945
$ws.Range("M9").Value = 54.64271666666667
$ws.Range("N9").Value = 163.92815
$ws.Range("O9").Value = 0.192030546333187
$ws.Range("P9").Value = 0.192030546333187
$ws.Range("Q9").Value = 1.026445218344445
$ws.Range("R9").Value = 9.2380069651
$ws.Range("S9").Value = 0.07462976730499239
$ws.Range("T9").Value = 0.07462976730499239
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.01878466666666667
$ws.Range("H10").Value = 0.056354
$ws.Range("I10").Value = 0.3886348746594945
$ws.Range("J10").Value = 0.3886348746594945
$ws.Range("M10").Value = 1.069012
$ws.Range("N10").Value = 3.207036
$ws.Range("O10").Value = 0.00375682196858928
$ws.Range("P10").Value = 0.00375682196858928
$ws.Range("Q10").Value = 0.02008103408266667
$ws.Range("R10").Value = 0.180729306744
$ws.Range("S10").Value = 0.00146003203488073
$ws.Range("T10").Value = 0.00146003203488073
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.01878466666666667
$ws.Range("H11").Value = 0.056354
$ws.Range("I11").Value = 0.3886348746594945
$ws.Range("J11").Value = 0.3886348746594945
$ws.Range("M11").Value = 54.36641700000001
$ws.Range("N11").Value = 163.099251
$ws.Range("O11").Value = 0.1910595481987908
$ws.Range("P11").Value = 0.1910595481987908
$ws.Range("Q11").Value = 1.021255021206
$ws.Range("R11").Value = 9.191295190854001
$ws.Range("S11").Value = 0.07425240356673669
$ws.Range("T11").Value = 0.07425240356673669